# Update the "Price" (column D) and "Volume(1h)" (column E) values for the
# cryptos sheet to reflect the latest GitHub Actions refresh.
#
# Column D values are stored as plain text (not numbers) in this workbook,
# so each cell's number format is forced to Text ("@") before the new
# value is written; otherwise Excel would auto-convert the numeric-looking
# string into a real number and silently drop meaningful trailing zeros
# (e.g. "0.8700" -> 0.87).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddress, $newValue) {
    $rng = $ws.Range($cellAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
}

# Column D (Price) updates
Set-TextValue "D2"  "244.73"
Set-TextValue "D3"  "23.92"
Set-TextValue "D4"  "5.199"
Set-TextValue "D5"  "0.05733"
Set-TextValue "D6"  "6.489"
Set-TextValue "D8"  "0.8145"
Set-TextValue "D9"  "0.8700"
Set-TextValue "D11" "0.06928"
Set-TextValue "D12" "0.03187"
Set-TextValue "D13" "0.02921"
Set-TextValue "D14" "0.09325"
Set-TextValue "D15" "3.853"
Set-TextValue "D16" "0.001530"
Set-TextValue "D17" "0.04715"
Set-TextValue "D18" "0.0005972"
Set-TextValue "D19" "0.006165"
Set-TextValue "D20" "0.001242"
Set-TextValue "D21" "0.004105"
Set-TextValue "D22" "0.00008498"
Set-TextValue "D24" "2.156"
Set-TextValue "D25" "0.3193"
Set-TextValue "D41" "0.006308"
Set-TextValue "D42" "0.1052"
Set-TextValue "D43" "0.002223"
Set-TextValue "D44" "0.008109"
Set-TextValue "D45" "0.00005470"
Set-TextValue "D48" "0.002564"

# Column E (Volume(1h)) updates - the "Bestin24h" / "Worstin24h" tag moved
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("E48").Value = "47BOLOBOLOBestin24h"
